$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 1750
$ws.Range("I52").Value = 1000
$ws.Range("J52").Value = 2000
$ws.Range("K52").Value = 3000
$ws.Range("L52").Value = 6000
$ws.Range("M52").Value = -2840
$ws.Range("N52").Value = -6320

$ws.Range("H132").Value = 14839430
$ws.Range("I132").Value = 1703.4706
$ws.Range("J132").Value = 183000340
$ws.Range("K132").Value = 5110.4118
$ws.Range("L132").Value = 549001020
$ws.Range("M132").Value = -2580.4118
$ws.Range("N132").Value = -549006080

$ws.Range("H137").Value = 3032671.2
$ws.Range("I137").Value = 6252244
$ws.Range("K137").Value = 18756732
$ws.Range("M137").Value = -18754182

$ws.Range("H138").Value = 1975666.1
$ws.Range("I138").Value = 48904.715
$ws.Range("J138").Value = 2607884.8
$ws.Range("K138").Value = 146714.145
$ws.Range("L138").Value = 7823654.399999999
$ws.Range("M138").Value = -141574.145
$ws.Range("N138").Value = -7833934.399999999

$ws.Range("H141").Value = 3257.1052
$ws.Range("I141").Value = 2876.4285
$ws.Range("K141").Value = 8629.2855
$ws.Range("M141").Value = -3449.2855


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 91092230
$ws.Range("I61").Value = 125126130
$ws.Range("J61").Value = 335171.34
$ws.Range("K61").Value = 125126130
$ws.Range("L61").Value = 335171.34
$ws.Range("M61").Value = -125125918
$ws.Range("N61").Value = -335595.34

$ws.Range("H122").Value = 2566.0588
$ws.Range("I122").Value = 2294.077
$ws.Range("J122").Value = 3450
$ws.Range("K122").Value = 6882.231000000001
$ws.Range("L122").Value = 10350
$ws.Range("M122").Value = -4432.231000000001
$ws.Range("N122").Value = -15250

$ws.Range("H132").Value = 47018.863
$ws.Range("I132").Value = 30121.885
$ws.Range("J132").Value = 112729.336
$ws.Range("K132").Value = 90365.655
$ws.Range("L132").Value = 338188.008
$ws.Range("M132").Value = -87835.655
$ws.Range("N132").Value = -343248.008

$ws.Range("H136").Value = 91092230
$ws.Range("I136").Value = 125126130
$ws.Range("J136").Value = 335171.34
$ws.Range("K136").Value = 375378390
$ws.Range("L136").Value = 1005514.02
$ws.Range("M136").Value = -375375840
$ws.Range("N136").Value = -1010614.02


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 10607.6
$ws.Range("I8").Value = 759.5
$ws.Range("J8").Value = 50000
$ws.Range("K8").Value = 759.5
$ws.Range("L8").Value = 50000
$ws.Range("M8").Value = -619.5
$ws.Range("N8").Value = -50280

$ws.Range("H134").Value = 4430.9062
$ws.Range("I134").Value = 3010.0476
$ws.Range("J134").Value = 7143.4546
$ws.Range("K134").Value = 9030.1428
$ws.Range("L134").Value = 21430.3638
$ws.Range("M134").Value = -6495.1428
$ws.Range("N134").Value = -26500.3638


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 30003162
$ws.Range("I58").Value = 35175156
$ws.Range("J58").Value = 5602.8
$ws.Range("K58").Value = 35175156
$ws.Range("L58").Value = 5602.8
$ws.Range("M58").Value = -35174953
$ws.Range("N58").Value = -6008.8

$ws.Range("H62").Value = 10666.667
$ws.Range("J62").Value = 10666.667
$ws.Range("L62").Value = 10666.667
$ws.Range("N62").Value = -11914.667

$ws.Range("H65").Value = 10666.667
$ws.Range("J65").Value = 10666.667
$ws.Range("L65").Value = 53333.335
$ws.Range("N65").Value = -59573.335

$ws.Range("H136").Value = 30003162
$ws.Range("I136").Value = 35175156
$ws.Range("J136").Value = 5602.8
$ws.Range("K136").Value = 105525468
$ws.Range("L136").Value = 16808.4
$ws.Range("M136").Value = -105522918
$ws.Range("N136").Value = -21908.4


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 811.6667
$ws.Range("I34").Value = 435
$ws.Range("K34").Value = 1305
$ws.Range("M34").Value = -1221

$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").ClearContents()

$ws.Range("H68").Value = 911.58826
$ws.Range("I68").Value = 564.5714
$ws.Range("J68").Value = 1001.55554
$ws.Range("K68").Value = 1693.7142
$ws.Range("L68").Value = 3004.66662
$ws.Range("M68").Value = -882.7142000000001
$ws.Range("N68").Value = -4626.66662

$ws.Range("H71").Value = 911.58826
$ws.Range("I71").Value = 564.5714
$ws.Range("J71").Value = 1001.55554
$ws.Range("K71").Value = 5081.1426
$ws.Range("L71").Value = 9013.99986
$ws.Range("M71").Value = -1025.1426
$ws.Range("N71").Value = -17125.99986

$ws.Range("H123").Value = 3265
$ws.Range("J123").Value = 3500
$ws.Range("L123").Value = 10500
$ws.Range("N123").Value = -15400

$ws.Range("H129").Value = 2253780.5
$ws.Range("I129").Value = 664.125
$ws.Range("J129").Value = 6413379.5
$ws.Range("K129").Value = 1992.375
$ws.Range("L129").Value = 19240138.5
$ws.Range("M129").Value = 3007.625
$ws.Range("N129").Value = -19250138.5

$ws.Range("H131").Value = 928.125
$ws.Range("J131").Value = 972.9
$ws.Range("L131").Value = 2918.7
$ws.Range("N131").Value = -12998.7

$ws.Range("H133").Value = 4362.273
$ws.Range("I133").Value = 3284.2856
$ws.Range("J133").Value = 6248.75
$ws.Range("K133").Value = 9852.856800000001
$ws.Range("L133").Value = 18746.25
$ws.Range("M133").Value = -4792.856800000001
$ws.Range("N133").Value = -28866.25


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 17201.334
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 25302
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 25302
$ws.Range("M3").Value = -884
$ws.Range("N3").Value = -25534

$ws.Range("H102").Value = 2372.5
$ws.Range("I102").Value = 2568.5
$ws.Range("J102").Value = 2078.5
$ws.Range("K102").Value = 2568.5
$ws.Range("L102").Value = 2078.5
$ws.Range("M102").Value = -946.5
$ws.Range("N102").Value = -5322.5


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7267.609
$ws.Range("I40").Value = 7772.5835
$ws.Range("J40").Value = 6716.727
$ws.Range("K40").Value = 7772.5835
$ws.Range("L40").Value = 6716.727
$ws.Range("M40").Value = -7636.5835
$ws.Range("N40").Value = -6988.727

$ws.Range("H136").Value = 71019.05
$ws.Range("I136").Value = 36268.965
$ws.Range("J136").Value = 148538.47
$ws.Range("K136").Value = 108806.895
$ws.Range("L136").Value = 445615.41
$ws.Range("M136").Value = -106256.895
$ws.Range("N136").Value = -450715.41


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H136").Value = 65556.19500000001
$ws.Range("I136").Value = 48602
$ws.Range("J136").Value = 101160
$ws.Range("K136").Value = 145806
$ws.Range("L136").Value = 303480
$ws.Range("M136").Value = -143256
$ws.Range("N136").Value = -308580

